$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlTop = -4160

# --- Row 10: only B/C content changes, formatting untouched ---
$ws.Range("B10").Value = 'Esta disciplina visa apresentar fundamentos de termodinâmica aplicada à área de ciência e engenharia de materiais. Especial ênfase é dada à energia na forma e calor para aquecimento de sistemas termodinâmicos; cálculos de variação de entalpia; entropia e energia de Gibbs de elementos e compostos em mudanças de estado; cálculos de variação de entalpia; entropia e energia de Gibbs de reação; aplicação da propriedade  energia de Gibbs para avaliação de transformações espontâneas e em equilíbrio; fundamentos de termodinâmica de soluções; cálculos de condições de equilíbrio em sistemas heterogêneos. Apresenta-se também as principais diferenças entre esta disciplina e a disciplina de Termodinâmica de Máquinas.'
$ws.Range("C10").Value = 'Esta disciplina visa apresentar fundamentos de termodinâmica aplicada à área de ciência e engenharia de materiais. Especial ênfase é dada à energia na forma e calor para aquecimento de sistemas termodinâmicos; cálculos de variação de entalpia; entropia e energia de Gibbs de elementos e compostos em mudanças de estado; cálculos de variação de entalpia; entropia e energia de Gibbs de reação; aplicação da propriedade  energia de Gibbs para avaliação de transformações espontâneas e em equilíbrio; fundamentos de termodinâmica de soluções; cálculos de condições de equilíbrio em sistemas heterogêneos. Apresenta-se também as principais diferenças entre esta disciplina e a disciplina de Termodinâmica de Máquinas.'

# --- Rows 13-26: content reshuffled, rebuild entirely ---
$ws.Range("A13:C24").ClearContents()

# Row 13
$ws.Range("B13").Value = '3577649 - Carlos Angelo Nunes'
$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = $xlTop
$ws.Range("C13").Value = '3577649 - Carlos Angelo Nunes'
$ws.Range("C13").Font.Bold = $false
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = $xlTop
$ws.Range("C13").Font.Color = 255
$ws.Rows.Item(13).AutoFit()

# Row 14
$ws.Range("B14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range("B14").Font.Bold = $false
$ws.Range("B14").WrapText = $true
$ws.Range("B14").VerticalAlignment = $xlTop
$ws.Range("C14").Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range("C14").Font.Bold = $false
$ws.Range("C14").WrapText = $true
$ws.Range("C14").VerticalAlignment = $xlTop
$ws.Range("C14").Font.Color = 255
$ws.Rows.Item(14).AutoFit()

# Row 15
$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("A15").Font.Bold = $true
$ws.Range("A15").WrapText = $false
$ws.Range("A15").VerticalAlignment = $xlTop
$ws.Range("B15").Value = '1) Introdução; 2) 1a Lei da Termodinâmica 3) 2a e 3a Leis da Termodinâmica 4) Equilíbrio heterogêneo: composição variável da fase gasosa; 5) Equilíbrio heterogêneo: composição variável da fase condensada;'
$ws.Range("B15").Font.Bold = $false
$ws.Range("B15").WrapText = $true
$ws.Range("B15").VerticalAlignment = $xlTop
$ws.Range("C15").Value = '1) Introdução; 2) 1a Lei da Termodinâmica 3) 2a e 3a Leis da Termodinâmica 4) Equilíbrio heterogêneo: composição variável da fase gasosa; 5) Equilíbrio heterogêneo: composição variável da fase condensada;'
$ws.Range("C15").Font.Bold = $false
$ws.Range("C15").WrapText = $true
$ws.Range("C15").VerticalAlignment = $xlTop
$ws.Range("C15").Font.Color = 255
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(15).RowHeight = 60

# Row 16
$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").WrapText = $false
$ws.Range("A16").VerticalAlignment = $xlTop
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(16).RowHeight = 60

# Row 17
$ws.Range("A17").Value = 'Programa:'
$ws.Range("A17").Font.Bold = $true
$ws.Range("A17").WrapText = $false
$ws.Range("A17").VerticalAlignment = $xlTop
$ws.Range("B17").Value = '1- Introdução: sistema; vizinhanças; fases; equilíbrio; fronteiras adiabáticas e diatérmicas; processos reversíveis e irreversíveis; estado termodinâmico; mudança de estado; processos cíclicos; equação de estado; calor; trabalho.2- A 1ª lei de Termodinâmica: energia interna; capacidades térmicas; entalpia; entalpia de transformação de fases; entalpia de formação e de reação; entalpia de reação em função da temperatura (introdução ao loop termodinâmico).3- A 2ª e 3ª leis da Termodinâmica: Dispersão de energia e entropia; entropia no zero absoluto; entropia de reação; entropia de reação em função da temperatura; desigualdade de Clausius; critérios de espontaneidade e equilíbrio; energia de Gibbs; energia de Helmholtz; energia de Gibbs de reação em função da temperatura; equação de Gibbs-Helmholtz.4- Equilíbrio heterogêneo: composição variável da fase gasosa: mistura de gases ideais; lei de Dalton; energia de Gibbs de um gás ideal; pressão de equilíbrio em sistemas metal-óxido-O2(g).5- Equilíbrio heterogêneo: composição variável da fase condensada: fugacidade; atividade termodinâmica; soluções e grandezas parciais molares; potencial químico; modelos de soluções; propriedades termodinâmicas de excesso'
$ws.Range("B17").Font.Bold = $false
$ws.Range("B17").WrapText = $true
$ws.Range("B17").VerticalAlignment = $xlTop
$ws.Range("C17").Value = '1- Introdução: sistema; vizinhanças; fases; equilíbrio; fronteiras adiabáticas e diatérmicas; processos reversíveis e irreversíveis; estado termodinâmico; mudança de estado; processos cíclicos; equação de estado; calor; trabalho.2- A 1ª lei de Termodinâmica: energia interna; capacidades térmicas; entalpia; entalpia de transformação de fases; entalpia de formação e de reação; entalpia de reação em função da temperatura (introdução ao loop termodinâmico).3- A 2ª e 3ª leis da Termodinâmica: Dispersão de energia e entropia; entropia no zero absoluto; entropia de reação; entropia de reação em função da temperatura; desigualdade de Clausius; critérios de espontaneidade e equilíbrio; energia de Gibbs; energia de Helmholtz; energia de Gibbs de reação em função da temperatura; equação de Gibbs-Helmholtz.4- Equilíbrio heterogêneo: composição variável da fase gasosa: mistura de gases ideais; lei de Dalton; energia de Gibbs de um gás ideal; pressão de equilíbrio em sistemas metal-óxido-O2(g).5- Equilíbrio heterogêneo: composição variável da fase condensada: fugacidade; atividade termodinâmica; soluções e grandezas parciais molares; potencial químico; modelos de soluções; propriedades termodinâmicas de excesso'
$ws.Range("C17").Font.Bold = $false
$ws.Range("C17").WrapText = $true
$ws.Range("C17").VerticalAlignment = $xlTop
$ws.Range("C17").Font.Color = 255
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("A18").Font.Bold = $true
$ws.Range("A18").WrapText = $false
$ws.Range("A18").VerticalAlignment = $xlTop
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(18).RowHeight = 120

# Row 19
$ws.Range("A19").Value = 'Avaliação:'
$ws.Range("A19").Font.Bold = $true
$ws.Range("A19").WrapText = $false
$ws.Range("A19").VerticalAlignment = $xlTop
$ws.Rows.Item(19).AutoFit()

# Row 20
$ws.Range("A20").Value = 'Método:'
$ws.Range("A20").Font.Bold = $true
$ws.Range("A20").WrapText = $false
$ws.Range("A20").VerticalAlignment = $xlTop
$ws.Range("B20").Value = 'Esta é uma disciplina fundamental, exigindo dedicação individual para assimilação de definições e conceitos. Isto envolve leitura concentrada e realização de exercícios numéricos.'
$ws.Range("B20").Font.Bold = $false
$ws.Range("B20").WrapText = $true
$ws.Range("B20").VerticalAlignment = $xlTop
$ws.Range("C20").Value = 'Esta é uma disciplina fundamental, exigindo dedicação individual para assimilação de definições e conceitos. Isto envolve leitura concentrada e realização de exercícios numéricos.'
$ws.Range("C20").Font.Bold = $false
$ws.Range("C20").WrapText = $true
$ws.Range("C20").VerticalAlignment = $xlTop
$ws.Range("C20").Font.Color = 255
$ws.Rows.Item(20).AutoFit()
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range("A21").Value = 'Critério:'
$ws.Range("A21").Font.Bold = $true
$ws.Range("A21").WrapText = $false
$ws.Range("A21").VerticalAlignment = $xlTop
$ws.Range("B21").Value = 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3'
$ws.Range("B21").Font.Bold = $false
$ws.Range("B21").WrapText = $true
$ws.Range("B21").VerticalAlignment = $xlTop
$ws.Range("C21").Value = 'Serão aplicadas duas avaliações escritas (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão:NF=(P1+2*P2)/3'
$ws.Range("C21").Font.Bold = $false
$ws.Range("C21").WrapText = $true
$ws.Range("C21").VerticalAlignment = $xlTop
$ws.Range("C21").Font.Color = 255
$ws.Rows.Item(21).AutoFit()
$ws.Rows.Item(21).RowHeight = 60

# Row 22
$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("A22").Font.Bold = $true
$ws.Range("A22").WrapText = $false
$ws.Range("A22").VerticalAlignment = $xlTop
$ws.Range("B22").Value = 'Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("B22").Font.Bold = $false
$ws.Range("B22").WrapText = $true
$ws.Range("B22").VerticalAlignment = $xlTop
$ws.Range("C22").Value = 'Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("C22").Font.Bold = $false
$ws.Range("C22").WrapText = $true
$ws.Range("C22").VerticalAlignment = $xlTop
$ws.Range("C22").Font.Color = 255
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(22).RowHeight = 60

# Row 23
$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("A23").Font.Bold = $true
$ws.Range("A23").WrapText = $false
$ws.Range("A23").VerticalAlignment = $xlTop
$ws.Range("B23").Value = '1) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995. ISBN 0-87339-270-1.2) P. Atkins & J. de Paula. Físico-Química, Livros Técnicos e Científicos Editora S.A., 2008. ISBN 978-85-216-1600-9.3) S.Stolen, T.Grande. Chemical Thermodynamics of Materials, John Wiley & Sons, Ltd. 2005. ISBN 978-0-471-49230-6.4) R. DeHoff. Thermodynamics in Materials Science. Taylor & Francis Group, 2006. ISBN 978-0-8493-4065-9.5) Y.A. Chang & W.A. Oates. Materials Thermodynamics, John Wiley & Sons, 2010. ISBN 978-0-470-48414-2.'
$ws.Range("B23").Font.Bold = $false
$ws.Range("B23").WrapText = $true
$ws.Range("B23").VerticalAlignment = $xlTop
$ws.Range("C23").Value = '1) Johnson, D.L. & Stracher, G.B., Thermodynamic - Loop Applications in Materials Systems, vols.1 e 2, The Minerals, Metals & Materials Society, 1995. ISBN 0-87339-270-1.2) P. Atkins & J. de Paula. Físico-Química, Livros Técnicos e Científicos Editora S.A., 2008. ISBN 978-85-216-1600-9.3) S.Stolen, T.Grande. Chemical Thermodynamics of Materials, John Wiley & Sons, Ltd. 2005. ISBN 978-0-471-49230-6.4) R. DeHoff. Thermodynamics in Materials Science. Taylor & Francis Group, 2006. ISBN 978-0-8493-4065-9.5) Y.A. Chang & W.A. Oates. Materials Thermodynamics, John Wiley & Sons, 2010. ISBN 978-0-470-48414-2.'
$ws.Range("C23").Font.Bold = $false
$ws.Range("C23").WrapText = $true
$ws.Range("C23").VerticalAlignment = $xlTop
$ws.Range("C23").Font.Color = 255
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(23).RowHeight = 120

# Row 24
$ws.Range("A24").Value = 'Requisitos:'
$ws.Range("A24").Font.Bold = $true
$ws.Range("A24").WrapText = $false
$ws.Range("A24").VerticalAlignment = $xlTop
$ws.Rows.Item(24).AutoFit()

# Row 25
$ws.Range("B25").Value = 'LOB1004 -  Cálculo II  (Requisito fraco)
'
$ws.Range("B25").Font.Bold = $false
$ws.Range("B25").WrapText = $true
$ws.Range("B25").VerticalAlignment = $xlTop
$ws.Range("C25").Value = 'LOB1004 -  Cálculo II  (Requisito fraco)
'
$ws.Range("C25").Font.Bold = $false
$ws.Range("C25").WrapText = $true
$ws.Range("C25").VerticalAlignment = $xlTop
$ws.Range("C25").Font.Color = 255
$ws.Rows.Item(25).AutoFit()
$ws.Rows.Item(25).RowHeight = 30

# Row 26
$ws.Range("B26").Value = 'LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'
$ws.Range("B26").Font.Bold = $false
$ws.Range("B26").WrapText = $true
$ws.Range("B26").VerticalAlignment = $xlTop
$ws.Range("C26").Value = 'LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito fraco)
'
$ws.Range("C26").Font.Bold = $false
$ws.Range("C26").WrapText = $true
$ws.Range("C26").VerticalAlignment = $xlTop
$ws.Range("C26").Font.Color = 255
$ws.Rows.Item(26).AutoFit()
$ws.Rows.Item(26).RowHeight = 30

